$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (prices in column D, 1h volume % in column E).
# Two coin rows (EthereumClassic/Toncoin and RenderToken/VeChain) also
# swapped rank order, so their B/C/D/E cells are rewritten in place.
#
# Column D holds plain-text price strings (e.g. "277.22", "1.000") in the
# source workbook. Excel auto-converts such numeric-looking text to real
# numbers on assignment, so we briefly force a text number format before
# writing the value and then clear the format again (these cells carry no
# custom formatting originally) to keep the cell text-typed without left
# -over style changes.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.806.20'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +4.18%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.875.24'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +3.21%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '277.22'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5289'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.92%  '

$ws.Range("E8").Value = '  -3.10%  '

$ws.Range("E9").Value = '  +4.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.08'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.11%  '

$ws.Range("E11").Value = '  -2.83%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07736'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.98%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.834.36'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.186'
$ws.Range("D14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.29'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +3.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.56'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +3.40%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008049'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9998'
$ws.Range("D19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.844.87'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.076.09'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.29%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.748'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.11%  '

$ws.Range("E23").Value = '  +0.38%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.188'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.52%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.383'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +7.97%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '146.55'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.21%  '

$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.664'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.45%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.34'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.35%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.62'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.77%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.347'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.43%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.309'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.86%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08907'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.60%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04926'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.91%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.173'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.36%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7275'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.877'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.31%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.281'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.77%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01859'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.35%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.327'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5141'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.68%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9518'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.36%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '116.12'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +5.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.171'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.80%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.114'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.00%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9996'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4481'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.52%  '

$ws.Range("E47").Value = '  -1.77%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.318'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.77%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.32'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05939'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.63%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.492'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.53%  '
